# Updates the cryptos list with latest prices / 1h volume changes.
# Some values (e.g. "601.08") look numeric, so Excel would silently
# convert them to real numbers on assignment. To preserve them as
# plain text (matching the original inline-string cells) we briefly
# force a text number-format before assigning, then restore the
# cell's style back to Normal/General so there is no visible/format
# side effect.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2"  "65.402.85"
Set-TextValue "E2"  "  -0.08%  "

Set-TextValue "D3"  "3.576.17"
Set-TextValue "E3"  "  +0.55%  "

Set-TextValue "E4"  "  -0.03%  "

Set-TextValue "D5"  "601.08"
Set-TextValue "E5"  "  -0.04%  "

Set-TextValue "D6"  "135.36"
Set-TextValue "E6"  "  -3.62%  "

Set-TextValue "D7"  "3.575.56"
Set-TextValue "E7"  "  +0.61%  "

Set-TextValue "E8"  "  +0.05%  "

Set-TextValue "E9"  "  +0.30%  "

Set-TextValue "E10" "  -1.20%  "

Set-TextValue "D11" "7.17"
Set-TextValue "E11" "  +1.90%  "

Set-TextValue "E12" "  -0.93%  "

Set-TextValue "D13" "4.188.14"
Set-TextValue "E13" "  +0.63%  "

Set-TextValue "D14" "0.0000184"
Set-TextValue "E14" "  -1.59%  "

Set-TextValue "D15" "27.53"
Set-TextValue "E15" "  +1.46%  "

Set-TextValue "D16" "3.577.66"
Set-TextValue "E16" "  +0.41%  "

Set-TextValue "E17" "  -0.05%  "

Set-TextValue "D18" "65.495.29"
Set-TextValue "E18" "  -0.04%  "

Set-TextValue "D19" "10.10"
Set-TextValue "E19" "  -1.72%  "

Set-TextValue "D20" "14.53"
Set-TextValue "E20" "  +1.81%  "

Set-TextValue "E21" "  -0.30%  "

Set-TextValue "D22" "392.73"
Set-TextValue "E22" "  -0.96%  "

Set-TextValue "D23" "0.583"
Set-TextValue "E23" "  +1.72%  "

Set-TextValue "D24" "3.721.08"
Set-TextValue "E24" "  +0.45%  "

Set-TextValue "D25" "74.23"
Set-TextValue "E25" "  -0.15%  "

Set-TextValue "E26" "  +0.06%  "

Set-TextValue "E27" "  -1.52%  "

Set-TextValue "D28" "8.09"
Set-TextValue "E28" "  +2.33%  "

Set-TextValue "D29" "1.67"
Set-TextValue "E29" "  +31.44%  "

Set-TextValue "D30" "8.65"
Set-TextValue "E30" "  +3.83%  "

# Rows 31 and 32 swap coin identity (Binance-PegBSC-USD <-> PancakeSwap)
Set-TextValue "B31" "PancakeSwap"
Set-TextValue "C31" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "2.31"
Set-TextValue "E31" "  +1.16%  "

Set-TextValue "B32" "Binance-PegBSC-USD"
Set-TextValue "C32" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D32" "0.999"
Set-TextValue "E32" "  -0.16%  "

Set-TextValue "D33" "3.584.19"
Set-TextValue "E33" "  +0.30%  "

Set-TextValue "D34" "24.30"
Set-TextValue "E34" "  +1.65%  "

Set-TextValue "E35" "  +0.02%  "

Set-TextValue "D36" "0.148"
Set-TextValue "E36" "  +0.00%  "

Set-TextValue "D37" "172.17"
Set-TextValue "E37" "  +2.56%  "

Set-TextValue "E38" "  -1.27%  "

Set-TextValue "E39" "  +2.48%  "

Set-TextValue "E40" "  +0.85%  "

Set-TextValue "D41" "0.0829"
Set-TextValue "E41" "  +2.91%  "

Set-TextValue "D42" "0.830"
Set-TextValue "E42" "  -0.37%  "

Set-TextValue "D43" "26.43"
Set-TextValue "E43" "  -1.03%  "

Set-TextValue "D44" "1.25"
Set-TextValue "E44" "  +4.93%  "

Set-TextValue "D45" "43.13"
Set-TextValue "E45" "  +0.31%  "

Set-TextValue "E46" "  -0.02%  "

Set-TextValue "E47" "  +0.77%  "

Set-TextValue "E48" "  -0.99%  "

Set-TextValue "E49" "  +2.19%  "

Set-TextValue "D50" "2.459.51"
Set-TextValue "E50" "  +0.33%  "

Set-TextValue "E51" "  +1.55%  "
